$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the second table (currently columns E:H) one column to the
# right, to F:I, leaving column E blank as a gap between the two tables.
# Move right-to-left so the source/destination ranges never overlap mid-move.
$ws.Range("H3:H46").Cut($ws.Range("I3"))
$ws.Range("G3:G46").Cut($ws.Range("H3"))
$ws.Range("F3:F46").Cut($ws.Range("G3"))
$ws.Range("E3:E46").Cut($ws.Range("F3"))

# --- The custom "20"-wide column setting (old column F) now belongs to the
# new column G.
$ws.Range("F1").ColumnWidth = $ws.Range("G1").ColumnWidth
$ws.Columns("G:G").ColumnWidth = 20

# --- Re-home the title. Stash its formatting on a scratch cell, strip the
# old anchor down to nothing, swap the merge over to the new (shifted left,
# smaller) area, then restore the formatting + text onto the new anchor.
$ws.Range("D1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D1").ClearFormats()

$ws.Range("D1:J2").UnMerge()
$ws.Range("C1:G2").Merge()

$ws.Range("Z1").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Z1").Clear()

$ws.Range("C1").Value = "Ginos Paysheet:04/2019"
$ws.Range("D1").ClearContents()
